# Add new columns I ("I0") and J ("IF") to Sheet1, matching the header
# style already used by the other header cells (B1:H1), and fill in the
# per-row numeric values for rows 2-59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell H1 (bold, centered,
# top-aligned, thin border) onto the two new header cells so they share
# the same style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data (rows 2-59), taken from the diff.
$iVals = @(5,9,11,7,9,8,9,8,7,9,7,8,7,9,9,7,5,8,9,8,8,9,9,8,9,7,7,9,8,8,9,8,7,8,6,7,7,7,8,9,8,7,8,8,9,8,8,10,10,9,7,6,7,6,6,5,4,3)
$jVals = @(6,9,11,7,10,9,9,8,8,9,8,8,7,9,9,7,6,8,9,8,8,9,9,9,9,8,7,9,8,8,9,8,8,8,7,8,7,7,8,9,8,7,8,8,9,9,8,11,10,10,7,6,7,6,6,5,4,3)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
